# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar las tasas de conversión en la nota de A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 9.89 = 40642.83 pesos", "1000 Bs = 9.63 = 39518.3 pesos")
$text = $text.Replace("40642.83 pesos = 9.86 = 956.88 Bs", "39518.3 pesos = 9.6 = 943.67 Bs")
$cellA1.Value() = $text

# --- tasas: actualizar las celdas numéricas de tasas ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value() = 103.8
$ws2.Range("O10").Value() = 4102
$ws2.Range("N12").Value() = 4117.9
$ws2.Range("O12").Value() = 98.333
